# "set up hybrid tests" - update MSRP staging data:
#  - bump a handful of 2020 model-year rows to 2021 with new base MSRP figures
#  - tweak several existing 2021 base MSRP figures
#  - append a new trim row (LC 500 Inspiration Series / 9260LE) at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 49: 9140 / LS 500h RWD -> now a 2021 model, new MSRP ---
$ws.Range("C49").Value = 2021
$ws.Range("D49").Value = 90500

# --- Row 50: 9146 / LS 500h AWD -> now a 2021 model, new MSRP ---
$ws.Range("C50").Value = 2021
$ws.Range("D50").Value = 93750

# --- Row 54: 9226 / RC F FUJI SPEEDWAY EDITION -> MSRP tweak ---
$ws.Range("D54").Value = 97625

# --- Row 80: LC Convertible Inspiration Series -> MSRP tweak ---
$ws.Range("D80").Value = 119900

# --- Rows 92-98: F SPORT Black Line Special Edition MSRP tweaks ---
$ws.Range("D92").Value = 49335
$ws.Range("D93").Value = 50735
$ws.Range("D94").Value = 51985
$ws.Range("D95").Value = 48835
$ws.Range("D96").Value = 51010
$ws.Range("D97").Value = 51765
$ws.Range("D98").Value = 53230

# --- New row 100: LC 500 Inspiration Series (9260LE) ---
# Set B before A so the shared-string table picks up "LC 500 INSPIRATION
# SERIES" before "9260LE", matching their relative index order.
$ws.Range("B100").Value = "LC 500 INSPIRATION SERIES"
$ws.Range("A100").Value = "9260LE"
$ws.Range("C100").Value = 2021
$ws.Range("D100").Value = 110420
$ws.Range("D100").NumberFormat = $ws.Range("D80").NumberFormat
$ws.Range("E100").Value = 1025
$ws.Range("E100").NumberFormat = $ws.Range("E80").NumberFormat

# --- View state: active cell / scroll position after the edits ---
$ws.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
